$d = $word.ActiveDocument

# Locate the two target paragraphs by their (unique) text content instead
# of a hard-coded index, so the script is resilient to small shifts.
$p1 = $null
$p2 = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*chaque période continue*") { $p1 = $p }
    if ($t -like "*La réservation peut s*") { $p2 = $p }
}

# ------------------------------------------------------------------
# Edit 1: "De la même manière, ... chaque période continue. "
#         -> drop the trailing space at the end of the sentence.
# ------------------------------------------------------------------
$r1 = $p1.Range
$trailingSpace = $d.Range($r1.End - 2, $r1.End - 1)
if ($trailingSpace.Text -eq " ") {
    $trailingSpace.Text = ""
}
# Touch (and immediately revert) formatting on the remaining text so the
# run keeps an explicit (empty) run-properties element after the edit.
$p1Body = $d.Range($r1.Start, $r1.End - 1)
$p1Body.Font.Bold = 1
$p1Body.Font.Bold = 0

# ------------------------------------------------------------------
# Edit 2: "La réservation peut s'effectuer jusqu'à un an à partir de
#          la date du jour de saisie. "
#         -> "Les réservations sont ouvertes sur une période de 365
#             jours."
# ------------------------------------------------------------------
$r2 = $p2.Range
$start = $r2.Start
$oldEnd = $r2.End - 1   # exclude the paragraph mark

# Wipe the old sentence.
$d.Range($start, $oldEnd).Text = ""

# Rebuild the new sentence as four runs:
#   "L" | "es" | " réservation" | "s sont ouvertes sur une période de 365 jours."
$d.Range($start, $start).InsertAfter("L")
$d.Range($start + 1, $start + 1).InsertAfter("es")
$d.Range($start + 3, $start + 3).InsertAfter(" réservation")
$d.Range($start + 15, $start + 15).InsertAfter("s sont ouvertes sur une période de 365 jours.")

$newEnd = $start + 15 + 46

# Nudge formatting on/off at each boundary so the four runs are kept
# distinct instead of being merged back together.
$d.Range($start, $start + 1).Font.Bold = 1
$d.Range($start, $start + 1).Font.Bold = 0

$d.Range($start + 1, $start + 3).Font.Bold = 1
$d.Range($start + 1, $start + 3).Font.Bold = 0

$d.Range($start + 3, $start + 15).Font.Bold = 1
$d.Range($start + 3, $start + 15).Font.Bold = 0

$d.Range($start + 15, $newEnd).Font.Bold = 1
$d.Range($start + 15, $newEnd).Font.Bold = 0
